$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("first")

$data = @(
    @(0.2352347881312385, 0.7011352911275793, 0.754111145231435,  0.868395730776836,  0.8601628886204244, 18),
    @(0.2672786791899959, 0.45307994943011,   0.2955295970660748, 0.5436263395624561, 0.4879522889190182, 17),
    @(0.3111056731803871, 0.3624139282624727, 0.1723564222806536, 0.4151583099019621, 0.2839148837627348, 16),
    @(0.3455294579606332, 0.3658231501549597, 0.1878025140189348, 0.4333618742101511, 0.2707369064906118, 15),
    @(0.3584323639489301, 0.3839945786029612, 0.2002730478501564, 0.4475187681540925, 0.2780688952452576, 14),
    @(0.364459496889083,  0.4106453173185426, 0.2107924313882408, 0.4591213689083103, 0.2906174726197602, 13),
    @(0.4143658927697875, 0.4310987738402203, 0.221430687247661,  0.4705642222350324, 0.2329219787203619, 12),
    @(0.3869142407525709, 0.4001865847611734, 0.1924388048854656, 0.4386784755210422, 0.2168174179053887, 11),
    @(0.3617893501639269, 0.3734161284103805, 0.1744359538854561, 0.4176553051087178, 0.2199606530304394, 10),
    @(0.3719529320823554, 0.3812107115717386, 0.1830055333568115, 0.4277914601260893, 0.2241397295910946, 9)
)

$rows = $data.Count
$cols = 6
$arr = New-Object 'object[,]' $rows,$cols

for ($i = 0; $i -lt $rows; $i++) {
    for ($j = 0; $j -lt $cols; $j++) {
        $arr[$i,$j] = $data[$i][$j]
    }
}

$ws.Range("B2:G11").Value = $arr
